# Updates cryptos list values (price, link, volume) per the authoritative diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$savedStyle = $ws.Range("D2").Style
$ws.Range("D2:E2").NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "41.797.92"
$ws.Cells.Item(2, 5).Value = "  +0.87%  "
$ws.Range("D2:E2").Style = $savedStyle

# Row 3
$savedStyle = $ws.Range("D3").Style
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.271.33"
$ws.Cells.Item(3, 5).Value = "  +0.74%  "
$ws.Range("D3:E3").Style = $savedStyle

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.03%  "

# Row 5
$savedStyle = $ws.Range("D5").Style
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "304.02"
$ws.Cells.Item(5, 5).Value = "  +0.47%  "
$ws.Range("D5:E5").Style = $savedStyle

# Row 6
$savedStyle = $ws.Range("D6").Style
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "92.66"
$ws.Cells.Item(6, 5).Value = "  +1.22%  "
$ws.Range("D6:E6").Style = $savedStyle

# Row 7
$ws.Cells.Item(7, 5).Value = "  +1.92%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.04%  "

# Row 9
$savedStyle = $ws.Range("D9").Style
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.485"
$ws.Cells.Item(9, 5).Value = "  -0.06%  "
$ws.Range("D9:E9").Style = $savedStyle

# Row 10
$savedStyle = $ws.Range("D10").Style
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "32.61"
$ws.Cells.Item(10, 5).Value = "  +1.63%  "
$ws.Range("D10:E10").Style = $savedStyle

# Row 11
$savedStyle = $ws.Range("D11").Style
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "53.57"
$ws.Cells.Item(11, 5).Value = "  -0.67%  "
$ws.Range("D11:E11").Style = $savedStyle

# Row 13
$ws.Cells.Item(13, 5).Value = "  -1.24%  "

# Row 14
$savedStyle = $ws.Range("D14").Style
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.70"
$ws.Cells.Item(14, 5).Value = "  +1.54%  "
$ws.Range("D14:E14").Style = $savedStyle

# Row 15
$savedStyle = $ws.Range("D15").Style
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "2.623.56"
$ws.Cells.Item(15, 5).Value = "  +0.82%  "
$ws.Range("D15:E15").Style = $savedStyle

# Row 16
$savedStyle = $ws.Range("D16").Style
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "14.29"
$ws.Cells.Item(16, 5).Value = "  +0.96%  "
$ws.Range("D16:E16").Style = $savedStyle

# Row 17
$savedStyle = $ws.Range("D17").Style
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "2.311.91"
$ws.Cells.Item(17, 5).Value = "  +5.67%  "
$ws.Range("D17:E17").Style = $savedStyle

# Row 18
$ws.Cells.Item(18, 5).Value = "  +3.56%  "

# Row 19
$savedStyle = $ws.Range("D19").Style
$ws.Range("D19:E19").NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "41.716.84"
$ws.Cells.Item(19, 5).Value = "  +0.91%  "
$ws.Range("D19:E19").Style = $savedStyle

# Row 20
$savedStyle = $ws.Range("D20").Style
$ws.Range("D20:E20").NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "12.55"
$ws.Cells.Item(20, 5).Value = "  +2.66%  "
$ws.Range("D20:E20").Style = $savedStyle

# Row 21
$savedStyle = $ws.Range("D21").Style
$ws.Range("D21:E21").NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "0.0₃0906"
$ws.Cells.Item(21, 5).Value = "  +0.17%  "
$ws.Range("D21:E21").Style = $savedStyle

# Row 22
$savedStyle = $ws.Range("D22").Style
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "5.93"
$ws.Cells.Item(22, 5).Value = "  +0.48%  "
$ws.Range("D22:E22").Style = $savedStyle

# Row 23
$ws.Cells.Item(23, 5).Value = "  +0.44%  "

# Row 24
$savedStyle = $ws.Range("D24").Style
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "243.78"
$ws.Cells.Item(24, 5).Value = "  +1.15%  "
$ws.Range("D24:E24").Style = $savedStyle

# Row 25
$ws.Cells.Item(25, 5).Value = "  +0.49%  "

# Row 26
$savedStyle = $ws.Range("D26").Style
$ws.Range("D26:E26").NumberFormat = "@"
$ws.Cells.Item(26, 2).Value = "ImmutableX"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(26, 4).Value = "1.94"
$ws.Cells.Item(26, 5).Value = "  +4.01%  "
$ws.Range("D26:E26").Style = $savedStyle

# Row 27
$savedStyle = $ws.Range("D27").Style
$ws.Range("D27:E27").NumberFormat = "@"
$ws.Cells.Item(27, 2).Value = "Dai"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(27, 4).Value = "1.00"
$ws.Cells.Item(27, 5).Value = "  -0.03%  "
$ws.Range("D27:E27").Style = $savedStyle

# Row 28
$savedStyle = $ws.Range("D28").Style
$ws.Range("D28:E28").NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "23.99"
$ws.Cells.Item(28, 5).Value = "  +0.92%  "
$ws.Range("D28:E28").Style = $savedStyle

# Row 29
$savedStyle = $ws.Range("D29").Style
$ws.Range("D29:E29").NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "9.48"
$ws.Cells.Item(29, 5).Value = "  -1.71%  "
$ws.Range("D29:E29").Style = $savedStyle

# Row 30
$savedStyle = $ws.Range("D30").Style
$ws.Range("D30:E30").NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "2.07"
$ws.Cells.Item(30, 5).Value = "  -5.38%  "
$ws.Range("D30:E30").Style = $savedStyle

# Row 31
$savedStyle = $ws.Range("D31").Style
$ws.Range("D31:E31").NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "35.39"
$ws.Cells.Item(31, 5).Value = "  +5.12%  "
$ws.Range("D31:E31").Style = $savedStyle

# Row 32
$savedStyle = $ws.Range("D32").Style
$ws.Range("D32:E32").NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "160.57"
$ws.Cells.Item(32, 5).Value = "  +1.61%  "
$ws.Range("D32:E32").Style = $savedStyle

# Row 33
$ws.Cells.Item(33, 5).Value = "  +0.96%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  -0.06%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +0.84%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -1.04%  "

# Row 37
$savedStyle = $ws.Range("D37").Style
$ws.Range("D37:E37").NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "16.91"
$ws.Cells.Item(37, 5).Value = "  +0.83%  "
$ws.Range("D37:E37").Style = $savedStyle

# Row 38
$ws.Cells.Item(38, 5).Value = "  +0.09%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +1.58%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  +0.77%  "

# Row 41
$savedStyle = $ws.Range("D41").Style
$ws.Range("D41:E41").NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.80"
$ws.Cells.Item(41, 5).Value = "  +1.17%  "
$ws.Range("D41:E41").Style = $savedStyle

# Row 42
$ws.Cells.Item(42, 5).Value = "  -1.50%  "

# Row 43
$savedStyle = $ws.Range("D43").Style
$ws.Range("D43:E43").NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.004.69"
$ws.Cells.Item(43, 5).Value = "  -2.93%  "
$ws.Range("D43:E43").Style = $savedStyle

# Row 44
$savedStyle = $ws.Range("D44").Style
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "19.49"
$ws.Cells.Item(44, 5).Value = "  -3.86%  "
$ws.Range("D44:E44").Style = $savedStyle

# Row 45
$savedStyle = $ws.Range("D45").Style
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.0281"
$ws.Cells.Item(45, 5).Value = "  +1.56%  "
$ws.Range("D45:E45").Style = $savedStyle

# Row 46
$savedStyle = $ws.Range("D46").Style
$ws.Range("D46:E46").NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "10.29"
$ws.Cells.Item(46, 5).Value = "  +1.16%  "
$ws.Range("D46:E46").Style = $savedStyle

# Row 47
$ws.Cells.Item(47, 5).Value = "  +3.00%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  -2.16%  "

# Row 49
$savedStyle = $ws.Range("D49").Style
$ws.Range("D49:E49").NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "52.68"
$ws.Cells.Item(49, 5).Value = "  +3.04%  "
$ws.Range("D49:E49").Style = $savedStyle

# Row 50
$savedStyle = $ws.Range("D50").Style
$ws.Range("D50:E50").NumberFormat = "@"
$ws.Cells.Item(50, 2).Value = "TrustWalletToken"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(50, 4).Value = "1.15"
$ws.Cells.Item(50, 5).Value = "  +1.03%  "
$ws.Range("D50:E50").Style = $savedStyle

# Row 51
$savedStyle = $ws.Range("D51").Style
$ws.Range("D51:E51").NumberFormat = "@"
$ws.Cells.Item(51, 2).Value = "Stacks"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(51, 4).Value = "1.51"
$ws.Cells.Item(51, 5).Value = "  -0.70%  "
$ws.Range("D51:E51").Style = $savedStyle
